$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder two pairs of countries (alphabetical re-sort of the list) ---
# Kirguistan now ranks above Afganistan (rows 55/56).
$ws.Range("A56").Value = "Afganistan"
$ws.Range("A55").Value = "Kirguistan"

# Uzbekistan now ranks above Marruecos (rows 64/65).
$ws.Range("A65").Value = "Marruecos"
$ws.Range("A64").Value = "Uzbekistan"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 08:13"

# --- Update statistic figures (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected
#     countries/rows ---

# Row 6 - India
$ws.Range("B6").Value = 1754117
$ws.Range("C6").Value = 2198
$ws.Range("D6").Value = 1148103
$ws.Range("E6").Value = 568599
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 37415

# Row 30 - Kazajistan
$ws.Range("B30").Value = 91593
$ws.Range("C30").Value = 1226
$ws.Range("D30").Value = 61839
$ws.Range("E30").Value = 28961
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 793

# Row 36 - Israel
$ws.Range("B36").Value = 72283
$ws.Range("C36").Value = 65
$ws.Range("D36").Value = 45629
$ws.Range("E36").Value = 26127
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 527

# Row 55 - now Kirguistan
$ws.Range("B55").Value = 36719
$ws.Range("C55").Value = 420
$ws.Range("D55").Value = 27274
$ws.Range("E55").Value = 8036
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 12
$ws.Range("H55").Value = 1409

# Row 56 - now Afganistan
$ws.Range("B56").Value = 36710
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 25509
$ws.Range("E56").Value = 9918
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 1283

# Row 64 - now Uzbekistan
$ws.Range("B64").Value = 25040
$ws.Range("C64").Value = 257
$ws.Range("D64").Value = 15299
$ws.Range("E64").Value = 9592
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 149

# Row 65 - now Marruecos
$ws.Range("B65").Value = 25015
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 17960
$ws.Range("E65").Value = 6688
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 367

# Row 74 - El Salvador
$ws.Range("D74").Value = 8561
$ws.Range("E74").Value = 8022
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 8
$ws.Range("H74").Value = 467

# Row 113 - Tailandia
$ws.Range("B113").Value = 3317
$ws.Range("C113").Value = 5
$ws.Range("D113").Value = 3142
$ws.Range("E113").Value = 117

Write-Output "edits applied"
